# Apply the "updated 4.0 files and mdl" edits to the Maximum Capacity Factor workbook.

$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the last-updated date in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45392

# --- "MCF" sheet: raise several capacity-factor inputs from 0.85/0.95 to 1 ---
$wsMcf = $wb.Worksheets.Item("MCF")
$wsMcf.Range("B2").Value = 1
$wsMcf.Range("B3").Value = 1
$wsMcf.Range("B4").Value = 1
$wsMcf.Range("B6").Value = 1
$wsMcf.Range("B10").Value = 1
$wsMcf.Range("B13").Value = 1
$wsMcf.Range("B14").Value = 1
$wsMcf.Range("B16").Value = 1
$wsMcf.Range("B17").Value = 1
$wsMcf.Range("B18").Value = 1

# --- update the active selection on the MCF sheet to B17 ---
$wsMcf.Activate() | Out-Null
$wsMcf.Range("B17").Select() | Out-Null
